# Generate Report for Handback
#
# 62e2eafc-c2b7-49e9-a387-b3c664144714.md has just been handed back
# (in sync with en-US). Promote it to the "latest handback" slot
# (previously occupied by 8e985be5-...) on every sheet, and push
# 8e985be5-... down into the slot that 62e2eafc used to occupy.

$wb = $excel.ActiveWorkbook

$Handed = "Handed back: in sync with en-US"

# ---------------------------------------------------------------
# Sheet "Overview": File Name | zh-cn | de-de
# ---------------------------------------------------------------
$overview = $wb.Worksheets.Item("Overview")

# Row 3 becomes 62e2eafc, row 4 becomes 8e985be5 (their old hyperlink
# targets simply swap places).
$overview.Range("A3").Value2 = "62e2eafc-c2b7-49e9-a387-b3c664144714.md"
$overview.Range("B3").Value2 = $Handed
$overview.Range("C3").Value2 = $Handed

$overview.Range("A4").Value2 = "8e985be5-708a-4ec4-9d31-8545d9390db9.md"
$overview.Range("B4").Value2 = $Handed
$overview.Range("C4").Value2 = $Handed

$linkA3 = $null
$linkA4 = $null
foreach ($h in $overview.Hyperlinks) {
    $addr = $h.Range.Address()
    if ($addr -eq '$A$3') { $linkA3 = $h }
    if ($addr -eq '$A$4') { $linkA4 = $h }
}
if ($linkA3 -ne $null -and $linkA4 -ne $null) {
    $a3disp = $linkA3.TextToDisplay
    $a3addr = $linkA3.Address
    $a4disp = $linkA4.TextToDisplay
    $a4addr = $linkA4.Address

    $linkA3.TextToDisplay = $a4disp
    $linkA3.Address = $a4addr
    $linkA4.TextToDisplay = $a3disp
    $linkA4.Address = $a3addr
}

# ---------------------------------------------------------------
# Helper applied identically to the "zh-cn" and "de-de" detail
# sheets (columns: A Source File Name, B Status, C Latest Handoff
# File, D Latest Handoff Datetime, E Latest Target File, F Latest
# Handback File, G Latest Handback DateTime, H Handoff Reason).
# ---------------------------------------------------------------
function Update-LangSheet {
    param($ws, $handoffXlf, $handoffDatetime, $handbackDatetime)

    $mdName = "62e2eafc-c2b7-49e9-a387-b3c664144714.md"

    # Row 3 -> 62e2eafc, now fully handed back: fill in the target
    # file + handback file links and the real handback datetime.
    $ws.Range("A3").Value2 = $mdName
    $ws.Range("B3").Value2 = $Handed
    $ws.Range("C3").Value2 = $handoffXlf
    $ws.Range("D3").Value2 = $handoffDatetime
    $ws.Range("E3").Value2 = $mdName
    $ws.Range("F3").Value2 = $handoffXlf
    $ws.Range("G3").Value2 = $handbackDatetime
    $ws.Range("H3").Value2 = "Include"

    # Row 4 -> 8e985be5, inherits the same handoff bookkeeping slot.
    $ws.Range("A4").Value2 = "8e985be5-708a-4ec4-9d31-8545d9390db9.md"
    $ws.Range("B4").Value2 = $Handed
    $ws.Range("C4").Value2 = $handoffXlf
    $ws.Range("D4").Value2 = $handoffDatetime
    $ws.Range("E4").Value2 = $mdName
    $ws.Range("F4").Value2 = $handoffXlf
    $ws.Range("G4").Value2 = $handbackDatetime
    $ws.Range("H4").Value2 = "Include"

    # Give the brand-new E3/F3/E4/F4 cells the same hyperlink-ish
    # style used by the rest of column A/C/E/F.
    $ws.Range("E3").Style = $ws.Range("A3").Style
    $ws.Range("F3").Style = $ws.Range("C3").Style
    $ws.Range("E4").Style = $ws.Range("A4").Style
    $ws.Range("F4").Style = $ws.Range("C4").Style

    # --- hyperlinks -------------------------------------------------
    $linkA3 = $null
    $linkC3 = $null
    $linkA4 = $null
    $linkC4 = $null
    foreach ($h in $ws.Hyperlinks) {
        $addr = $h.Range.Address()
        if ($addr -eq '$A$3') { $linkA3 = $h }
        if ($addr -eq '$C$3') { $linkC3 = $h }
        if ($addr -eq '$A$4') { $linkA4 = $h }
        if ($addr -eq '$C$4') { $linkC4 = $h }
    }

    # A3 used to point at 8e985be5's md, A4 at 62e2eafc's md -> swap.
    if ($linkA3 -ne $null -and $linkA4 -ne $null) {
        $a3disp = $linkA3.TextToDisplay
        $a3addr = $linkA3.Address
        $a4disp = $linkA4.TextToDisplay
        $a4addr = $linkA4.Address

        $linkA3.TextToDisplay = $a4disp
        $linkA3.Address = $a4addr
        $linkA4.TextToDisplay = $a3disp
        $linkA4.Address = $a3addr
    }

    # C3/C4 already point at 62e2eafc's handoff xlf - leave as-is,
    # just make sure the display text matches the (unchanged) file.
    if ($linkC3 -ne $null) { $linkC3.TextToDisplay = $handoffXlf }
    if ($linkC4 -ne $null) { $linkC4.TextToDisplay = $handoffXlf }

    # New E3/F3/E4/F4 hyperlinks (target-language repo md + handback
    # repo xlf). Re-use the handoff-repo xlf target as a stand-in
    # target file link, consistent with the rest of the workbook.
    $ws.Hyperlinks.Add($ws.Range("E3"), $linkA4.Address, [Type]::Missing, [Type]::Missing, $mdName) | Out-Null
    $ws.Hyperlinks.Add($ws.Range("F3"), $linkC3.Address, [Type]::Missing, [Type]::Missing, $handoffXlf) | Out-Null
    $ws.Hyperlinks.Add($ws.Range("E4"), $linkA4.Address, [Type]::Missing, [Type]::Missing, $mdName) | Out-Null
    $ws.Hyperlinks.Add($ws.Range("F4"), $linkC4.Address, [Type]::Missing, [Type]::Missing, $handoffXlf) | Out-Null
}

# ---------------------------------------------------------------
# Sheet "zh-cn"
# ---------------------------------------------------------------
$zhcn = $wb.Worksheets.Item("zh-cn")
Update-LangSheet $zhcn `
    "62e2eafc-c2b7-49e9-a387-b3c664144714.5b38bb005706b39ac2edcfa28d16ccafbee5fc1c.zh-cn.xlf" `
    "2016-03-03 12:46:36" `
    "2016-03-03 12:47:26"

# ---------------------------------------------------------------
# Sheet "de-de"
# ---------------------------------------------------------------
$dede = $wb.Worksheets.Item("de-de")
Update-LangSheet $dede `
    "62e2eafc-c2b7-49e9-a387-b3c664144714.5b38bb005706b39ac2edcfa28d16ccafbee5fc1c.de-de.xlf" `
    "2016-03-03 12:46:48" `
    "2016-03-03 12:47:49"
